$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "TextBox 24" ---
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Split the first paragraph ("Web Application") into several runs that
# together read "Web Application (Gradio, Streamlit, etc.)", keeping the
# bold formatting of the original run.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Web "
$run2 = $para1.InsertAfter("Application (")
$run3 = $run2.InsertAfter("Gradio")
$run4 = $run3.InsertAfter(", ")
$run5 = $run4.InsertAfter("Streamlit")
$run6 = $run5.InsertAfter(", etc.)")

# Reposition/resize the textbox (PowerPoint re-lays this out because of
# the extra wrapped line; match the canonical values exactly).
$shp.Top = 85.7623

# --- Shape 3: "Straight Arrow Connector 7" ---
# Its start is glued to TextBox 24's bottom connection point, so its
# length grows a little once TextBox 24 has been resized above.
$conn = $s.Shapes.Item(3)
$conn.Height = 5587 / 12700
